# PS2_tables.xlsx - "update PS2 to make nice for Latex"
# The "Sheet2" tab holds the New-College NPV/tax-revenue table. The
# "New College" row (row 11) previously discounted I11 (=F11/1.03) into the
# NPV total; that extra one-year discount term is removed, the number of
# new-college workers (M11) is corrected from 9.3 to 7 (thousands), and the
# stray "Total" label in N12 is replaced with an actual grand-total formula.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet2")
$ws.Activate()

# Remove the one-year-discount term that fed into the New College NPV (L11)
$ws.Range("I11").ClearContents()

# Correct the number of new college workers (in thousands)
$ws.Range("M11").Value = 7

# Turn the "Total" text label into a real grand total formula
$ws.Range("N12").Formula = "=N11+N9"

# Leave the selection where the author finished editing
$ws.Range("N12").Select() | Out-Null
$excel.ActiveWindow.Zoom = 100
